$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep header row (A1/B1) untouched, rewrite the data rows (2-8) with the
# updated list of Instagram usernames, all with "OK" status.
$names = @(
    "pearl.what.is.this.behaviour",
    "bcbilliofficial",
    "aman.gupta.09",
    "dhattarwalaman",
    "Benc4n",
    "tumblrindeed",
    "thesavagebean"
)

$row = 2
foreach ($name in $names) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = "OK"
    $row++
}

$ws.Range("B5:B8").Select()
